$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New time-log entry in row 105 (previously blank except the date-column style)
$ws.Range("A105").Value = 41945
$ws.Range("B105").Value = 0.0097222222222222224
$ws.Range("C105").Value = 0.1013888888888889
$ws.Range("D105").Value = 20
$ws.Range("E105").Formula = "=IF(AND(NOT(ISBLANK(B105)),NOT(ISBLANK(C105))),(C105-B105)*24-D105/60,"""")"
$ws.Range("F105").Value = "Coding"

# Move the active selection from C105 to C106, matching where editing continued
$ws.Range("C106").Select()
